# Applies the progress-report-group self-check prompt revisions.
# Uses Find/Replace (not wildcard) on the whole document story so that
# matches can span multiple runs.

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $range = $d.Content
    $range.Find.ClearFormatting()
    $ok = $range.Find.Execute($find, $true, $true, $false, $false, $false, `
                               $true, 1, $false, $replace, 2)
    if (-not $ok) {
        throw "Find/Replace failed for: $find"
    }
}

# 1. Opening paragraph: note the draft is about the Recommendation Report,
#    and tweak "check the" -> "check for the".
Replace-Text "Progress Report on the work our group has completed. I am uploading" `
             "Progress Report on the work our group has completed on our Recommendation Report. I am uploading"
Replace-Text "Please help me check the following criteria, but do NOT" `
             "Please help me check for the following criteria, but do NOT"

# 2. Drop the stray trailing space at the end of the memo-format paragraph.
Replace-Text "a signature). " "a signature)."

# 3. Introduction paragraph: drop the Fact Sheet Collection mention, drop
#    "at least", and clarify which paragraph covers which content.
Replace-Text "reminds the reader of the Recommendation Report or Fact Sheet Collection project." `
             "reminds the reader of the Recommendation Report project."
Replace-Text "The introduction is at least two paragraphs long." `
             "The introduction is two paragraphs long."
Replace-Text "summarizes the work that has been accomplished so far, and provides a purpose statement, which identifies the purpose of the document." `
             "summarizes the work that has been accomplished so far in the first paragraph, and provides a purpose statement, which identifies the purpose of the document, in the second paragraph."

# 4. Work Completed section: replace the task-summary guidance with
#    subheading-formatting guidance (moved up from the headings paragraph).
Replace-Text " The Work Completed section begins with a brief introduction that summarizes the work our group has completed so far. It then summarizes the main tasks completed to date, specifying the time period covered and the major tasks completed [ex: writing, research, production of visual aids]. This section then discusses each major task that has been completed in a paragraph or two of its own, with its own subheading. The completed major tasks in the order in which they were brought up in the introduction to the section." `
             " The section organizes the information with subheadings, using these guidelines: It includes a brief sentence or two of introduction between the main heading for the section and the first subheading. The subheadings are visibly different from the Work Completed section heading. They use a slightly smaller font than the main section heading (but a bigger font than is used for the paragraphs). They can also be a different color or size."

# 5. Work Scheduled section: simplify to Gantt-chart-first guidance.
Replace-Text "introduces the schedule with a brief summary of the work that still needs to be done, organizing the work chronologically and specifying the time period covered. The Work Scheduled Section then discusses each scheduled major task in a paragraph of its own, with its own subheading. The scheduled major tasks are organized in the order they were brought up in the introduction to the section. The section ends with a Gantt Chart that breaks down the remaining tasks and gives a tentative completion date for each." `
             "introduces the schedule with a simple explanation of the information to follow. This section includes a Gantt Chart with an updated schedule for the group and concludes with a brief summary of the work that still needs to be done."

# 6. Conclusion section: trim down to just the contact-information requirement.
Replace-Text " The Conclusion should do at least one of the following: appraise the work done thus far, draw conclusions about the work, or make recommendations concerning the work. If more than one of the options is used, the information is organized with separate subheadings. The Conclusion looks to future tasks in a sentence or two before ending, indicating flexibility and encouraging reader response. The Conclusion ends by providing specific" `
             " It provides specific"

# 7. Headings paragraph: drop "and subheadings" plus the subheading-formatting
#    guidance (now relocated to the Work Completed paragraph above).
Replace-Text "Our group progress report uses informative headings and subheadings to help your readers navigate the report. The subheadings are visibly different from the main section heading. They use a slightly smaller font than the main section heading (but a bigger font than is used for the paragraphs). They can also be a different color or size." `
             "Our group progress report uses informative headings to help your readers navigate the report."
